# Title Designs.pptx edit
# Commit: "3.4 and 3.5 renamed properly / Finished editing decanters besides git issues."
#
# Semantic change: two new "divider" slides are inserted into the
# "Chapter 3: Separators" section (right after the existing
# "3.1 Flash Separators" slide, i.e. at position 4):
#   - new slide "3.4 Decanters"
#   - new slide "3.5 Ethyl Acetate Plant " / "Distillation Column"
# All later slides shift down by two positions but are otherwise
# unchanged. The slide-master / slide-layout "last saved" date field
# is also bumped from 9/23/2021 to 9/28/2021.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Insert the two new divider slides by duplicating the existing
#    "Chapter 3: Separators / 3.1 Flash Separators" slide (slide 3),
#    which already has the right layout, shapes, and formatting.
#    Duplicate() inserts the copy immediately after its source, so
#    duplicating slide 3 twice in a row lands both new slides at
#    positions 4 and 5, pushing everything else down by two slots.
# ---------------------------------------------------------------
$dividerSlide = $p.Slides.Item(3)

$dup1Range = $dividerSlide.Duplicate()
$slide34 = $dup1Range.Item(1)

$dup2Range = $slide34.Duplicate()
$slide35 = $dup2Range.Item(1)

# ---------------------------------------------------------------
# 2. Slide "3.4 Decanters" (new position 4)
#    Shape 1 = Title ("Chapter 3: Separators" - already correct)
#    Shape 2 = subtitle TextBox -> needs its text replaced.
# ---------------------------------------------------------------
$slide34.Shapes.Item(2).TextFrame.TextRange.Text = "3.4 Decanters"

# ---------------------------------------------------------------
# 3. Slide "3.5 Ethyl Acetate Plant / Distillation Column"
#    (new position 5). The subtitle is split across two runs in the
#    source file ("3.5 Ethyl Acetate Plant " then "Distillation
#    Column"), matching the pattern used elsewhere in the deck
#    (e.g. "5.3 Ethyl Acetate Plant " / "Heat Exchanger").
# ---------------------------------------------------------------
$subtitle35 = $slide35.Shapes.Item(2).TextFrame.TextRange
$part1 = "3.5 Ethyl Acetate Plant "
$part2 = "Distillation Column"
$subtitle35.Text = $part1 + $part2

# Re-apply formatting to the second run explicitly so it becomes its
# own <a:r> run (matching the two-run structure in the target file)
# while keeping the same font / size / color as the rest of the line.
$run2 = $subtitle35.Characters($part1.Length + 1, $part2.Length)
$run2.Font.Name = "Palatino Linotype"
$run2.Font.Size = 32
$run2.Font.Color.RGB = $subtitle35.Characters(1,1).Font.Color.RGB

# ---------------------------------------------------------------
# 4. Bump the cached "last saved" date field (datetimeFigureOut)
#    shown on the slide master and every slide layout from
#    9/23/2021 to 9/28/2021.
# ---------------------------------------------------------------
$newDate = "9/28/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "9/23/2021") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "9/23/2021") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Write-Output "Slides now: $($p.Slides.Count)"
